$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, shifting existing rows 31-60 down to 32-61.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new record's data.
$ws.Cells.Item(31, 1).Value = 10
$ws.Cells.Item(31, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(31, 3).Value = "La Araucanía"
$ws.Cells.Item(31, 4).Value = 44413
$ws.Cells.Item(31, 5).Value = 9
$ws.Cells.Item(31, 6).Value = 100112031
$ws.Cells.Item(31, 7).Value = "Poroto verde"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 40
$ws.Cells.Item(31, 11).Value = 35000
$ws.Cells.Item(31, 12).Value = 35000
$ws.Cells.Item(31, 13).Value = 35000
$ws.Cells.Item(31, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 1400
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
